$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = "69.262.53"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(2, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +1.35%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.785.18"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(3, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.52%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(4, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.65%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = "628.24"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(5, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +4.44%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = "163.93"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(6, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.90%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.781.39"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(7, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.60%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(8, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.04%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.519"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(9, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.10%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.160"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(10, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.63%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(11, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.08%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(12, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +2.36%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(13, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -1.01%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = "35.53"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(14, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.97%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.440.89"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(15, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.03%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.765.33"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(16, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.76%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = "69.361.55"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(17, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +1.53%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = "17.95"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(18, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -2.72%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(19, 2)
$cell.NumberFormat = "@"
$cell.Value = "TRON"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(19, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.113"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(19, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -1.01%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(20, 2)
$cell.NumberFormat = "@"
$cell.Value = "Polkadot"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(20, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.07"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(20, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.39%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = "467.91"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(21, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +1.23%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = "9.61"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(22, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -1.12%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.703"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(23, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.14%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(24, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.93%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = "83.27"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(25, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.18%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = "12.03"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(26, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.36%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.14"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(27, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +1.56%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(28, 2)
$cell.NumberFormat = "@"
$cell.Value = "Dai"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(28, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(28, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.02%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(29, 2)
$cell.NumberFormat = "@"
$cell.Value = "RenderToken"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(29, 3)
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = "10.00"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(29, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.05%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.955.58"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(30, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.00%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.67"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(31, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.60%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.21"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(32, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.75%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.26"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(33, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -1.40%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = "28.93"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(34, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -1.58%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(35, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.17%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.737.56"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(36, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.47%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = "8.98"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(37, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.53%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(38, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +2.74%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(39, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +8.27%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.30"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(40, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.78%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(41, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.63%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.968"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(42, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -2.37%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(43, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.01%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(44, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.04%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = "155.06"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(45, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +2.19%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.298"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(46, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.86%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.93"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(47, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +2.59%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = "46.97"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(48, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -1.17%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = "42.56"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(49, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -1.38%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(50, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +2.69%  "
$cell.Style = "Normal"
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = "8.41"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(51, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.30%  "
$cell.Style = "Normal"
